# Capitalize the hex-byte letters (a-f -> A-F) in the "doip" (G) and "uds" (H)
# columns, e.g. "0x02:0xfd:0x00" -> "0x02:0xFD:0x00". The leading "0x" prefix
# stays lowercase; "N/A" values are left untouched.

function Uppercase-HexCodes($value) {
    if ($value -eq "N/A") {
        return $value
    }

    $bytes = $value -split ':'
    $upperBytes = @()
    foreach ($b in $bytes) {
        $digits = $b.Substring(2)
        $upperBytes += "0x" + $digits.ToUpper()
    }
    return ($upperBytes -join ':')
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 38; $row++) {
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = Uppercase-HexCodes $gCell.Value2

    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = Uppercase-HexCodes $hCell.Value2
}
